$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest scraped values.
# Cells hold plain text (not numbers/percent values), so NumberFormat is forced to
# Text ("@") right before the assignment to stop Excel from re-typing numeric-looking
# strings (e.g. "0.996") as actual numbers; Style is then reset to "Normal" so the
# cell keeps its original (default) appearance.

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '42.659.15'
$c.Style = "Normal"
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.345.76'
$c.Style = "Normal"
$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  -1.65%  '
$c.Style = "Normal"
$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  -0.07%  '
$c.Style = "Normal"
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '319.25'
$c.Style = "Normal"
$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  -0.90%  '
$c.Style = "Normal"
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '105.72'
$c.Style = "Normal"
$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  +0.43%  '
$c.Style = "Normal"
$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  -3.21%  '
$c.Style = "Normal"
$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  +0.06%  '
$c.Style = "Normal"
$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  -5.63%  '
$c.Style = "Normal"
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '41.18'
$c.Style = "Normal"
$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  -1.33%  '
$c.Style = "Normal"
$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  -2.08%  '
$c.Style = "Normal"
$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  -1.73%  '
$c.Style = "Normal"
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '0.996'
$c.Style = "Normal"
$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  -1.63%  '
$c.Style = "Normal"
$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  -0.23%  '
$c.Style = "Normal"
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '15.88'
$c.Style = "Normal"
$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  -7.82%  '
$c.Style = "Normal"
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '2.702.87'
$c.Style = "Normal"
$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  -1.64%  '
$c.Style = "Normal"
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '2.336.54'
$c.Style = "Normal"
$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  -2.60%  '
$c.Style = "Normal"
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '42.676.82'
$c.Style = "Normal"
$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  -1.33%  '
$c.Style = "Normal"
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '7.69'
$c.Style = "Normal"
$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  +4.52%  '
$c.Style = "Normal"
$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  -2.37%  '
$c.Style = "Normal"
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '77.07'
$c.Style = "Normal"
$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  +1.57%  '
$c.Style = "Normal"
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '3.62'
$c.Style = "Normal"
$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  +5.16%  '
$c.Style = "Normal"
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '257.63'
$c.Style = "Normal"
$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  -3.82%  '
$c.Style = "Normal"
$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  -5.01%  '
$c.Style = "Normal"
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '9.37'
$c.Style = "Normal"
$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  -4.36%  '
$c.Style = "Normal"
$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  +0.13%  '
$c.Style = "Normal"
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '11.34'
$c.Style = "Normal"
$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  -4.04%  '
$c.Style = "Normal"
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '22.89'
$c.Style = "Normal"
$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  +0.21%  '
$c.Style = "Normal"
$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  +0.73%  '
$c.Style = "Normal"
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '174.69'
$c.Style = "Normal"
$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  -1.22%  '
$c.Style = "Normal"
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '36.19'
$c.Style = "Normal"
$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -3.64%  '
$c.Style = "Normal"
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.0886'
$c.Style = "Normal"
$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -4.46%  '
$c.Style = "Normal"
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '6.08'
$c.Style = "Normal"
$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  +3.40%  '
$c.Style = "Normal"
$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  -8.10%  '
$c.Style = "Normal"
$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  +11.67%  '
$c.Style = "Normal"
$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  -3.40%  '
$c.Style = "Normal"
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '4.58'
$c.Style = "Normal"
$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -5.62%  '
$c.Style = "Normal"
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.0360'
$c.Style = "Normal"
$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  -2.19%  '
$c.Style = "Normal"
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '3.75'
$c.Style = "Normal"
$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  -8.29%  '
$c.Style = "Normal"
$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  -4.50%  '
$c.Style = "Normal"
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '71.66'
$c.Style = "Normal"
$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  +3.55%  '
$c.Style = "Normal"
$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  +0.61%  '
$c.Style = "Normal"
$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  -7.54%  '
$c.Style = "Normal"
$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  -0.16%  '
$c.Style = "Normal"
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '114.69'
$c.Style = "Normal"
$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  -8.95%  '
$c.Style = "Normal"
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '11.88'
$c.Style = "Normal"
$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  -4.84%  '
$c.Style = "Normal"
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '9.11'
$c.Style = "Normal"
$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  -5.09%  '
$c.Style = "Normal"
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '72.94'
$c.Style = "Normal"
$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c.Style = "Normal"
$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  -1.56%  '
$c.Style = "Normal"

# Ranking shuffled: THORChain moved up to row 47, BitcoinSV dropped to row 48
# (link/price/volume all move with their coin).
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '5.48'
$c.Style = "Normal"
$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -2.82%  '
$c.Style = "Normal"

$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '85.77'
$c.Style = "Normal"
$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  -0.70%  '
$c.Style = "Normal"
